# Add a new slide 19 ("Anas Zouhir - Contributions") to the deck.
#
# The new slide re-uses the exact shape layout/formatting of the last
# existing "Contributions" slide (slide 18 - "Camron Darpoh"), so the
# cleanest, most faithful way to reproduce it through the PowerPoint
# object model is to duplicate slide 18 (which places the duplicate
# immediately after it, i.e. as the new slide 19) and then replace the
# title and bullet text.

$p = $ppt.ActivePresentation

$sourceSlide = $p.Slides.Item($p.Slides.Count)
$newSlide = $sourceSlide.Duplicate().Item(1)

$titleShape = $null
$contentShape = $null
for ($i = 1; $i -le $newSlide.Shapes.Count; $i++) {
    $sh = $newSlide.Shapes.Item($i)
    if ($sh.Name -eq "Title 1") {
        $titleShape = $sh
    }
    if ($sh.Name -eq "Content Placeholder 2") {
        $contentShape = $sh
    }
}

# Fall back to the placeholder shapes by position if the duplicated
# slide's shape names ever differ from the expected ones.
if ($titleShape -eq $null -or $contentShape -eq $null) {
    for ($i = 1; $i -le $newSlide.Shapes.Count; $i++) {
        $sh = $newSlide.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($titleShape -eq $null -and $sh.Name -like "Title*") {
                $titleShape = $sh
            } elseif ($contentShape -eq $null) {
                $contentShape = $sh
            }
        }
    }
}

$titleShape.TextFrame.TextRange.Text = "Anas Zouhir - Contributions"

$bulletLines = @(
    "Created Front End of the App",
    "Designed Class Diagram for the system",
    "Coded the home page class (display timetable)",
    "Coded Map feature of the app",
    "Currently working on local database (Shared Preferences)"
)
$contentShape.TextFrame.TextRange.Text = [string]::Join([char]13, $bulletLines) + [char]13 + [char]13
